$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helpers -----------------------------------------------------------
# Excel's COM .Value setter auto-detects dates/numbers in strings like
# "01/01/2023" and silently converts the cell into a numeric date. To keep
# such values as literal text we stage them in a scratch cell with a
# leading apostrophe (forces text), copy just the *value* into the real
# destination, then copy *formatting* back in from a known-good sibling
# cell so the destination's style (s="2"/s="3") is preserved.
function Set-LiteralText($destCell, $text, $formatSourceCell) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $destCell.PasteSpecial(-4163) # xlPasteValues
    $formatSourceCell.Copy()
    $destCell.PasteSpecial(-4122) # xlPasteFormats
    $scratch.Clear()
}

# Cells that previously had no content at all come back from a plain
# ".Value = " assignment with the wrong (inherited) style, so re-apply the
# correct column style from a known-good sibling afterwards.
function Fix-Style($destCell, $formatSourceCell) {
    $formatSourceCell.Copy()
    $destCell.PasteSpecial(-4122) # xlPasteFormats
}

# --- "Name:" row -> new discipline name ---------------------------------
$ws.Range("B4").Value = "Graduation Monograph II"
$ws.Range("C4").Value = "Graduation Monograph II"

# --- "Ativação:" row -> updated activation date (keep as literal text) --
Set-LiteralText $ws.Range("B8") "01/01/2023" $ws.Range("B9")
Set-LiteralText $ws.Range("C8") "01/01/2023" $ws.Range("C9")

# --- "Objetivos:" row -> new responsible professor ----------------------
$ws.Range("B10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C10").Value = "5840730 - Antonio Jefferson da Silva Machado"

# --- "Objectives:" row -> new English objectives text (was empty) -------
$ws.Range("B11").Value = "The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer."
$ws.Range("C11").Value = "The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer."
Fix-Style $ws.Range("B11") $ws.Range("B10")

# --- "Programa resumido:" row -> now carries the activation date value --
Set-LiteralText $ws.Range("B13") "01/01/2023" $ws.Range("B9")
Set-LiteralText $ws.Range("C13") "01/01/2023" $ws.Range("C9")

# --- "Short syllabus:" row -> new English short syllabus text (was empty)
$ws.Range("B14").Value = "Prepare a monograph of Undergraduate Work under the guidance of a professor and present it to a panel of examiners."
$ws.Range("C14").Value = "Prepare a monograph of Undergraduate Work under the guidance of a professor and present it to a panel of examiners."
Fix-Style $ws.Range("B14") $ws.Range("B13")

# --- "Programa:" row -> now carries the responsible professor value -----
$ws.Range("B15").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C15").Value = "5840730 - Antonio Jefferson da Silva Machado"

# --- "Syllabus:" row -> new English syllabus text (was empty) -----------
$ws.Range("B16").Value = "The course program will consist of the following steps: 1) Preparation and writing of a monograph on a previously defined and approved subject in the Undergraduate Work I discipline. 2) Definition and disclosure of the presentation date after delivery of the monograph in advance of at least , 15 working days. 3) Definition of the panel of examiners, consisting of the supervisor and at least two invited professionals, with training in engineering or related areas. 4) Presentation and evaluation of the TG. 5) Publication of the evaluation. In case of approval, the final copy of the monograph (printed and electronic copy) must be delivered with the agreement of the supervisor."
$ws.Range("C16").Value = "The course program will consist of the following steps: 1) Preparation and writing of a monograph on a previously defined and approved subject in the Undergraduate Work I discipline. 2) Definition and disclosure of the presentation date after delivery of the monograph in advance of at least , 15 working days. 3) Definition of the panel of examiners, consisting of the supervisor and at least two invited professionals, with training in engineering or related areas. 4) Presentation and evaluation of the TG. 5) Publication of the evaluation. In case of approval, the final copy of the monograph (printed and electronic copy) must be delivered with the agreement of the supervisor."
Fix-Style $ws.Range("B16") $ws.Range("B15")

# --- "Método:" row -> new responsible professor --------------------------
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
